# Auto-generated edit script: updates numeric cells in Kraken_Profits workbook
# per the authoritative diff (scheduled runner refresh of profit calcs).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 26783.334
$ws.Range("I18").Value = 24140
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 24140
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = -23856
$ws.Range("N18").Value = -40568
# row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 2072.2856
$ws.Range("I41").Value = 3467.6667
$ws.Range("J41").Value = 1025.75
$ws.Range("K41").Value = 3467.6667
$ws.Range("L41").Value = 1025.75
$ws.Range("M41").Value = -3027.6667
$ws.Range("N41").Value = -1905.75
# row 124 (Leve Item ID 34241)
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815
# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 3205.682
$ws.Range("I132").Value = 3205.682
$ws.Range("K132").Value = 9617.045999999998
$ws.Range("M132").Value = -7087.045999999998
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3822.2222
$ws.Range("J138").Value = 3993.75
$ws.Range("L138").Value = 11981.25
$ws.Range("N138").Value = -22261.25

$ws = $wb.Worksheets.Item("ARM")
# row 36 (Leve Item ID 3068)
$ws.Range("H36").Value = 3155.3333
$ws.Range("I36").Value = 3155.3333
$ws.Range("K36").Value = 3155.3333
$ws.Range("M36").Value = -2809.3333
# row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 65000
$ws.Range("J37").Value = 65000
$ws.Range("L37").Value = 65000
$ws.Range("N37").Value = -65546
# row 76 (Leve Item ID 10679)
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# row 79 (Leve Item ID 10679)
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# row 125 (Leve Item ID 34251)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# row 127 (Leve Item ID 35347)
$ws.Range("H127").Value = 99995
$ws.Range("J127").Value = 99995
$ws.Range("L127").Value = 99995
$ws.Range("N127").Value = -109915

$ws = $wb.Worksheets.Item("BSM")
# row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 800
$ws.Range("I99").Value = 800
$ws.Range("K99").Value = 800
$ws.Range("M99").Value = 698
# row 127 (Leve Item ID 35358)
$ws.Range("H127").Value = 99995
$ws.Range("J127").Value = 99995
$ws.Range("L127").Value = 99995
$ws.Range("N127").Value = -109915

$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4717.125
$ws.Range("I31").Value = 1999.75
$ws.Range("J31").Value = 5622.9165
$ws.Range("K31").Value = 1999.75
$ws.Range("L31").Value = 5622.9165
$ws.Range("M31").Value = -1704.75
$ws.Range("N31").Value = -6212.9165
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4717.125
$ws.Range("I34").Value = 1999.75
$ws.Range("J34").Value = 5622.9165
$ws.Range("K34").Value = 1999.75
$ws.Range("L34").Value = 5622.9165
$ws.Range("M34").Value = -1797.75
$ws.Range("N34").Value = -6026.9165
# row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 9999
$ws.Range("I41").Value = 9999
$ws.Range("K41").Value = 9999
$ws.Range("M41").Value = -9571
# row 51 (Leve Item ID 2039)
$ws.Range("H51").Value = 24000
# row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 19624.75
$ws.Range("I60").Value = 16998
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 16998
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -16487
$ws.Range("N60").Value = -21022
# row 61 (Leve Item ID 2039)
$ws.Range("H61").Value = 24000
# row 80 (Leve Item ID 12015)
$ws.Range("H80").Value = 90000
$ws.Range("J80").Value = 90000
$ws.Range("L80").Value = 90000
$ws.Range("N80").Value = -92246
# row 83 (Leve Item ID 12015)
$ws.Range("H83").Value = 90000
$ws.Range("J83").Value = 90000
$ws.Range("L83").Value = 270000
$ws.Range("N83").Value = -281232

$ws = $wb.Worksheets.Item("CUL")
# row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 2155
$ws.Range("I129").Value = 1342.2858
$ws.Range("J129").Value = 4999.5
$ws.Range("K129").Value = 4026.8574
$ws.Range("L129").Value = 14998.5
$ws.Range("M129").Value = 973.1425999999997
$ws.Range("N129").Value = -24998.5
# row 130 (Leve Item ID 36058)
$ws.Range("H130").Value = 4000
$ws.Range("J130").Value = 3500
$ws.Range("L130").Value = 10500
$ws.Range("N130").Value = -20540
# row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1278
$ws.Range("J131").Value = 2350
$ws.Range("L131").Value = 7050
$ws.Range("N131").Value = -17130

$ws = $wb.Worksheets.Item("GSM")
# row 4 (Leve Item ID 2056)
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224
# row 11 (Leve Item ID 4422)
$ws.Range("H11").Value = 11751125
$ws.Range("I11").Value = 13287000
$ws.Range("J11").Value = 1000000
$ws.Range("K11").Value = 13287000
$ws.Range("L11").Value = 1000000
$ws.Range("M11").Value = -13286861
$ws.Range("N11").Value = -1000278
# row 12 (Leve Item ID 4093)
$ws.Range("H12").Value = 2749.75
$ws.Range("I12").Value = 2749.75
$ws.Range("K12").Value = 2749.75
$ws.Range("M12").Value = -2609.75
# row 14 (Leve Item ID 4198)
$ws.Range("H14").Value = 6601450
$ws.Range("I14").Value = 8501917
$ws.Range("J14").Value = 3750749.8
$ws.Range("K14").Value = 8501917
$ws.Range("L14").Value = 3750749.8
$ws.Range("M14").Value = -8501749
$ws.Range("N14").Value = -3751085.8
# row 15 (Leve Item ID 12018)
$ws.Range("H15").Value = 75000
$ws.Range("J15").Value = 75000
$ws.Range("L15").Value = 75000
$ws.Range("N15").Value = -75576
# row 27 (Leve Item ID 2061)
$ws.Range("H27").Value = 1199
$ws.Range("J27").Value = 1199
$ws.Range("L27").Value = 1199
$ws.Range("N27").Value = -1531
# row 81 (Leve Item ID 12018)
$ws.Range("H81").Value = 75000
$ws.Range("J81").Value = 75000
$ws.Range("L81").Value = 75000
$ws.Range("N81").Value = -76996
# row 84 (Leve Item ID 12018)
$ws.Range("H84").Value = 75000
$ws.Range("J84").Value = 75000
$ws.Range("L84").Value = 225000
$ws.Range("N84").Value = -234984
# row 95 (Leve Item ID 18235)
$ws.Range("H95").Value = 44288.8
$ws.Range("J95").Value = 44288.8
$ws.Range("L95").Value = 44288.8
$ws.Range("N95").Value = -49780.8
# row 99 (Leve Item ID 19532)
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 1246
# row 100 (Leve Item ID 18367)
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164
# row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 2389.6428
$ws.Range("I102").Value = 2389.6428
$ws.Range("K102").Value = 2389.6428
$ws.Range("M102").Value = -767.6428000000001

$ws = $wb.Worksheets.Item("LTW")
# row 11 (Leve Item ID 3542)
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -860
$ws.Range("N11").ClearContents()
# row 18 (Leve Item ID 3772)
$ws.Range("H18").Value = 15199.8
$ws.Range("I18").Value = 14000
$ws.Range("J18").Value = 19999
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 19999
$ws.Range("M18").Value = -13828
$ws.Range("N18").Value = -20343
# row 26 (Leve Item ID 3559)
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# row 69 (Leve Item ID 10671)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# row 72 (Leve Item ID 10671)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# row 88 (Leve Item ID 10961)
$ws.Range("H88").Value = 29000
$ws.Range("I88").Value = 29000
$ws.Range("K88").Value = 29000
$ws.Range("M88").Value = -28572
# row 91 (Leve Item ID 10961)
$ws.Range("H91").Value = 29000
$ws.Range("I91").Value = 29000
$ws.Range("K91").Value = 29000
$ws.Range("M91").Value = -27518

$ws = $wb.Worksheets.Item("WVR")
# row 14 (Leve Item ID 2658)
$ws.Range("H14").Value = 10001.333
$ws.Range("I14").Value = 10004
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 10004
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -9836
$ws.Range("N14").Value = -10336
# row 48 (Leve Item ID 3140)
$ws.Range("H48").Value = 10019.667
$ws.Range("I48").Value = 10029.5
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 10029.5
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = -9460.5
$ws.Range("N48").Value = -11138
# row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 54995
$ws.Range("I70").Value = 54990
$ws.Range("K70").Value = 54990
$ws.Range("M70").Value = -54675
# row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 54995
$ws.Range("I73").Value = 54990
$ws.Range("K73").Value = 54990
$ws.Range("M73").Value = -53898
# row 75 (Leve Item ID 11957)
$ws.Range("H75").Value = 70000
$ws.Range("I75").Value = 70000
$ws.Range("K75").Value = 70000
$ws.Range("M75").Value = -69064
# row 78 (Leve Item ID 11957)
$ws.Range("H78").Value = 70000
$ws.Range("I78").Value = 70000
$ws.Range("K78").Value = 210000
$ws.Range("M78").Value = -205320
# row 86 (Leve Item ID 11977)
$ws.Range("H86").Value = 80325
$ws.Range("J86").Value = 80325
$ws.Range("L86").Value = 80325
$ws.Range("N86").Value = -82571
# row 89 (Leve Item ID 11977)
$ws.Range("H89").Value = 80325
$ws.Range("J89").Value = 80325
$ws.Range("L89").Value = 401625
$ws.Range("N89").Value = -412857

